$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the document (end of the
#    last paragraph) to the very start of the document (before the first
#    run of the first paragraph). Bookmarks.Add() with a degenerate range
#    that sits exactly at document position 0 has a quirk where it snaps to
#    cover the whole first paragraph, so we work around it by temporarily
#    inserting a single placeholder character at position 0, anchoring the
#    bookmark to that character, and then deleting the character again -
#    which correctly collapses the bookmark back to a zero-length range at
#    position 0.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$p1.InsertBefore("X")
$startRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $startRange)
$startRange.Delete()

# ---------------------------------------------------------------------------
# 2) Extend the last paragraph ("... solución híbrida.") with the new
#    sentences, including a spell-check "Clouds" run wrapped in proofErr
#    markers, exactly as the target run layout requires. Because plain text
#    insertion (InsertAfter/TypeText) always coalesces into the neighboring
#    run when formatting matches, and InsertXML only ever inserts whole
#    paragraphs, we build the final paragraph (original text + new runs)
#    as one InsertXML call right after the existing paragraph, and then
#    delete the old paragraph's content+mark so only the newly built
#    paragraph remains in its place.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$oldPara = $d.Paragraphs($lastParaIndex).Range
Write-Output ("old paragraph [" + $oldPara.Start + "," + $oldPara.End + "]: [" + $oldPara.Text + "]")

$newParaXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005C45DE" w:rsidRDefault="005C45DE"><w:r><w:tab/><w:t>Por temas de flexibilidad, escalabilidad o seguridad y acceso, nos encontraremos que en ocasiones preferimos optar por una soluci' + [char]0x00F3 + 'n h' + [char]0x00ED + 'brida.</w:t></w:r><w:r><w:t xml:space="preserve"> Habr' + [char]0x00E1 + ' cierta necesidad de servicios o aplicaciones que queramos mantener bajo control total con infraestructuras locales. Luego aquellos servicios que puedan ser m' + [char]0x00E1 + 's susceptibles de cambios o modificaciones de recursos necesarios, es decir, que necesiten de adaptabilidad o escalabilidad, ya sea tanto de crecimiento o decrecimiento, son las que dejaremos en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Clouds</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> p' + [char]0x00FA + 'blicas, para no tener que abordar constantemente problemas con adquisici' + [char]0x00F3 + 'n de nuevo hardware o equipamiento. La flexibilidad de este tipo de combinaciones, nos permiten ser mucho m' + [char]0x00E1 + 's concretos a la hora de resolver cualquier problema con nuestras necesidades.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint = $d.Range($oldPara.End, $oldPara.End)
$insertPoint.InsertXML($newParaXml)

# Re-fetch the (now stale) old paragraph - still at the same index, right
# before our freshly inserted paragraph - and remove it completely
# (content + its own paragraph mark), leaving only the new paragraph behind.
$oldPara2 = $d.Paragraphs($lastParaIndex).Range
Write-Output ("deleting old paragraph [" + $oldPara2.Start + "," + $oldPara2.End + "]: [" + $oldPara2.Text + "]")
$oldPara2.Delete()

Write-Output "done"
